$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = -2
    4  = 1
    6  = 6
    8  = 1
    9  = -3
    10 = 3
    11 = -1
    12 = -3
    13 = -2
    14 = 7
    15 = -1
    16 = -1
    17 = -1
    18 = -4
    20 = -2
    21 = 2
    22 = 3
    23 = 1
    24 = 1
    26 = -4
    28 = -3
    29 = -2
    30 = -5
    31 = 1
    32 = 3
    34 = 1
    35 = 1
    37 = -1
    38 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
